# Refresh the cryptocurrency snapshot values (Price / Volume(1h)) in columns D and E.
# Price values that look like plain numbers are written with a leading apostrophe so
# Excel keeps them as literal text (matching the original inline-string formatting)
# instead of re-parsing "24.50" -> 24.5 or "0.00001363" -> scientific notation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.669.72"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.695.79"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("D5").Value = "'315.38"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").Value = "'0.3918"
$ws.Range("E7").Value = "  -0.87%  "
$ws.Range("D8").Value = "'0.4056"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "'1.502"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").Value = "'1.006"
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("D11").Value = "'52.98"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "'0.08767"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").Value = "'7.659"
$ws.Range("E13").Value = "  +5.20%  "
$ws.Range("D14").Value = "'24.50"
$ws.Range("E14").Value = "  +3.33%  "
$ws.Range("D15").Value = "'0.00001363"
$ws.Range("E15").Value = "  +3.12%  "
$ws.Range("D16").Value = "'7.993"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").Value = "1.696.24"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "'98.50"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("D19").Value = "'0.07115"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").Value = "'19.84"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").Value = "'7.381"
$ws.Range("E21").Value = "  +4.62%  "
$ws.Range("D22").Value = "'1.007"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").Value = "'14.32"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "24.671.22"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "'3.031"
$ws.Range("E25").Value = "  -7.03%  "
$ws.Range("D26").Value = "'2.357"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "'22.76"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").Value = "'162.85"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").Value = "'8.494"
$ws.Range("E29").Value = "  +12.92%  "
$ws.Range("D30").Value = "'137.57"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").Value = "'5.233"
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("D32").Value = "1.883.22"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "'0.08920"
$ws.Range("E33").Value = "  +3.64%  "
$ws.Range("D34").Value = "'7.534"
$ws.Range("E34").Value = "  +5.44%  "
$ws.Range("D35").Value = "'1.054"
$ws.Range("E36").Value = "  +3.76%  "
$ws.Range("D37").Value = "'0.02934"
$ws.Range("E37").Value = "  +7.42%  "
$ws.Range("D38").Value = "'0.2742"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").Value = "'10.81"
$ws.Range("E39").Value = "  -5.42%  "
$ws.Range("D40").Value = "'14.33"
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("D41").Value = "'0.09139"
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("D42").Value = "'0.7924"
$ws.Range("E42").Value = "  +3.08%  "
$ws.Range("D43").Value = "'1.470"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "'16.78"
$ws.Range("E44").Value = "  +4.11%  "
$ws.Range("D45").Value = "'0.7249"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").Value = "'2.579"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").Value = "'4.219"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "'1.005"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").Value = "'1.332"
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("D50").Value = "'139.38"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").Value = "'91.38"
$ws.Range("E51").Value = "  +1.01%  "
